$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("D1").Value = "price"
$ws.Range("F1").Value = "status"

# Row 2 (date in C2 is unchanged, leave as-is)
$ws.Range("A2").Value = "tirth"
$ws.Range("B2").Value = "Grindlays Super Saver Income Fund-GSSIF-Half Yearly Dividend"
$ws.Range("D2").Value = 15.5
$ws.Range("E2").Value = 500
$ws.Range("F2").Value = "Buy"

# Row 3 (date in C3 is unchanged, leave as-is)
$ws.Range("A3").Value = "yamik"
$ws.Range("B3").Value = "Aditya Birla Sun Life Gilt Plus - Liquid Plan - Growth - Regular Plan"
$ws.Range("D3").Value = 98.75
$ws.Range("E3").Value = 75000
$ws.Range("F3").Value = "Sale"

# Row 4 (date in C4 is unchanged, leave as-is)
$ws.Range("A4").Value = "xyz"
$ws.Range("B4").Value = "cdshszjkcg"
$ws.Range("D4").Value = 120.3
$ws.Range("E4").Value = 60000
$ws.Range("F4").Value = "Buy"

# Apply the header formatting (bold font, thin border, centered/top aligned) to F1
# to match the rest of row 1 (A1:E1 already use this formatting)
$f1 = $ws.Range("F1")
$f1.Borders.LineStyle = 1        # xlContinuous
$f1.Borders.Weight = 2           # xlThin
$f1.Font.Bold = $true
$f1.HorizontalAlignment = -4108  # xlCenter
$f1.VerticalAlignment = -4160    # xlTop

# Leave the selection on the newly added column header, matching the
# final cursor position left behind by the edit session
$f1.Select() | Out-Null
